# Applies the "Update team members: delete rank and modify David's
# affectation" edit to equipe-leptonex.pptx.
#
# Summary of changes:
#  1. Refresh the "datetimeFigureOut" date placeholder text (03/04/2025 ->
#     09/04/2025) on the slide master and every slide layout.
#  2. On slide 1, shape 1 ("Google Shape;157;p13", the academic-partners
#     list): strip the trailing " (RANK ORG)" suffix from each team
#     member's name, leaving just the plain name.
#  3. On slide 1, shape 2 ("Google Shape;158;p13", the non-academic
#     partners list): change David Gomis' affectation text and resize /
#     reposition the shape's bounding box.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update the date placeholder ("03/04/2025" -> "09/04/2025") on the
#    slide master and on every slide layout.
# ---------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "03/04/2025") {
                $shp.TextFrame.TextRange.Text = "09/04/2025"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes($master.Shapes)

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateShapes($layouts.Item($i).Shapes)
}

# ---------------------------------------------------------------------
# 2. Slide 1 / shape 1: drop the " (RANK ORG)" suffix from each name.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$namesShape = $slide.Shapes.Item(1)
$namesTr = $namesShape.TextFrame.TextRange

$nameParaIdx = @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 29)
foreach ($i in $nameParaIdx) {
    $para = $namesTr.Paragraphs($i, 1)
    $t = $para.Text
    $cut = $t.IndexOf(" (")
    if ($cut -ge 0) {
        $para.Text = $t.Substring(0, $cut)
    }
}

# ---------------------------------------------------------------------
# 3. Slide 1 / shape 2: update David Gomis' affectation text and move /
#    resize the shape.
# ---------------------------------------------------------------------
$orgShape = $slide.Shapes.Item(2)
$orgTr = $orgShape.TextFrame.TextRange

$davidPara = $orgTr.Paragraphs(7, 1)
$prefixLen = "`tDavid Gomis".Length
$oldSuffixLen = $davidPara.Length - 1 - $prefixLen
$suffixStart = $davidPara.Start + $prefixLen
$suffixRange = $orgTr.Characters($suffixStart, $oldSuffixLen)
$suffixRange.Text = ", Pôle Biodiversité Paysages Agroécologie et Alimentation"

$orgShape.Left = 451.8935039370079
$orgShape.Width = 476.39220472440945
